# Applies the changes described by the commit diff to the open presentation.
#
# Slide 4 ("Updates Since Version-04"), Content Placeholder 2:
#   - Paragraph "Discuss review comment on using IOAM FEC (SFL) for HbH IOAM"
#     becomes "Use GAL with different IOAM G-ACh for E2E and HbH"
#   - Paragraph "Discuss multiple G-ACh / Control Word handling"
#     becomes "Discuss multiple G-ACh / Control Word headers"
#   - A new paragraph "Discuss using IOAM FEC (SFL) for HbH IOAM" is added
#     right after it (same bullet style).
#
# Slide 6, "TextBox 2":
#   - "Next IP Version is added in GACH Header " becomes
#     "Next IP Version is added in G-ACh Header "
#   - "Next Protocol = 0000 or 0001 or IPV4 0100 or IPv6 0110, etc." becomes
#     "Next IP Version = 0000 or 0001 or IPV4 0100 or IPv6 0110, etc."

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 4 - "Open Items" bullets
# ---------------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$content = $slide4.Shapes.Item(2)
$tr = $content.TextFrame.TextRange

# --- Paragraph: "Discuss review comment on using IOAM FEC (SFL) for HbH IOAM"
$para9 = $tr.Paragraphs(9)
$para9.Text = "Use GAL with different IOAM G-ACh for E2E and HbH IOAM"

# Split out "ACh" as its own run.
$t = $para9.Text
$achStart = $t.IndexOf("ACh") + 1
$para9.Characters($achStart, 3).Text = "ACh"

# Merge the text between "ACh" and "HbH" into a single run.
$t = $para9.Text
$achEnd = $t.IndexOf("ACh") + 3
$hbhStart = $t.IndexOf("HbH")
$midStart = $achEnd + 1
$midLen = $hbhStart - $achEnd
$para9.Characters($midStart, $midLen).Text = " for E2E and "

# Drop the trailing " IOAM" that used to follow "HbH".
$t = $para9.Text
$tailStart = $t.LastIndexOf(" IOAM") + 1
$tailLen = $t.Length - $t.LastIndexOf(" IOAM")
$para9.Characters($tailStart, $tailLen).Delete()

# --- Paragraph: "Discuss multiple G-ACh / Control Word handling"
$para10 = $tr.Paragraphs(10)
$t = $para10.Text
$marker = " / Control Word handling"
$markerStart = $t.IndexOf($marker) + 1
$para10.Characters($markerStart, $marker.Length).Text = " / Control Word headers"

# --- New paragraph after it: "Discuss using IOAM FEC (SFL) for HbH IOAM"
$para10 = $tr.Paragraphs(10)
$para10.InsertAfter([char]13 + "Discuss using IOAM FEC (SFL) for HbH IOAM")
$para11 = $tr.Paragraphs(11)
$t = $para11.Text
$hbhStart = $t.IndexOf("HbH") + 1
$para11.Characters($hbhStart, 3).Text = "HbH"

# ---------------------------------------------------------------------------
# Slide 6 - GACH Header textbox
# ---------------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$textbox = $slide6.Shapes.Item(5)
$tr6 = $textbox.TextFrame.TextRange

# --- Paragraph: "Next IP Version is added in GACH Header "
$paraA = $tr6.Paragraphs(1)
$t = $paraA.Text
$gachStart = $t.IndexOf("GACH Header") + 1
$paraA.Characters($gachStart, "GACH".Length).Text = "G-ACh"

# Split "G-ACh" into "G-" and "ACh" runs.
$t = $paraA.Text
$achStart = $t.IndexOf("ACh") + 1
$paraA.Characters($achStart, 3).Text = "ACh"

# Merge the leading text with "G-" back into a single run.
$t = $paraA.Text
$achStart0 = $t.IndexOf("ACh")
$paraA.Characters(1, $achStart0).Text = $paraA.Characters(1, $achStart0).Text

# --- Paragraph: "Next Protocol = 0000 or 0001 or IPV4 0100 or IPv6 0110, etc."
$paraB = $tr6.Paragraphs(3)
$paraB.Text = "RESET_PLACEHOLDER_TEXT"
$paraB.Text = "Next IP Version = 0000 or 0001 or IPV4 0100 or IPv6 0110, etc."
